{"js": "const body = context.document.body;\nconst replacements = [\n  [\"2025-03-18 Tuesday\", \"2025-03-19 Wednesday\"],\n  [\"12\u00d774=\", \"12\u00d782=\"],\n  [\"49\u00d747=\", \"56\u00d757=\"],\n  [\"88\u00d794=\", \"46\u00d784=\"],\n  [\"88\u00d723=\", \"56\u00d750=\"],\n  [\"51\u00d735=\", \"47\u00d763=\"],\n  [\"38\u00d752=\", \"75\u00d786=\"],\n  [\"75\u00d792=\", \"89\u00d731=\"],\n  [\"49\u00d757=\", \"51\u00d742=\"],\n  [\"79\u00d740=\", \"63\u00d773=\"],\n  [\"60\u00d796=\", \"94\u00d759=\"],\n  [\"58\u00d744=\", \"39\u00d718=\"],\n  [\"45\u00d784=\", \"83\u00d790=\"],\n  [\"60\u00d754=\", \"92\u00d796=\"],\n  [\"85\u00d774=\", \"34\u00d763=\"],\n  [\"33\u00d780=\", \"22\u00d740=\"],\n  [\"19\u00d733=\", \"72\u00d718=\"],\n  [\"16\u00d763=\", \"50\u00d725=\"],\n  [\"70\u00d756=\", \"80\u00d776=\"],\n  [\"89\u00d739=\", \"16\u00d738=\"],\n  [\"64\u00d792=\", \"49\u00d723=\"],\n  [\"98\u00d718=\", \"61\u00d723=\"],\n  [\"20\u00d716=\", \"67\u00d740=\"],\n  [\"19\u00d762=\", \"34\u00d764=\"],\n  [\"91\u00d775=\", \"63\u00d787=\"],\n  [\"84\u00d777=\", \"32\u00d738=\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\nfunction Replace-Text($old, $new) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $old\n    $find.Replacement.Text = $new\n    $find.Execute([ref]$find.Text, $false, $false, $false, $false, $false, $true, 1, $false, [ref]$find.Replacement.Text, 2)\n}\n\nReplace-Text '2025-03-18 Tuesday' '2025-03-19 Wednesday'\nReplace-Text '12\u00d774=' '12\u00d782='\nReplace-Text '49\u00d747=' '56\u00d757='\nReplace-Text '88\u00d794=' '46\u00d784='\nReplace-Text '88\u00d723=' '56\u00d750='\nReplace-Text '51\u00d735=' '47\u00d763='\nReplace-Text '38\u00d752=' '75\u00d786='\nReplace-Text '75\u00d792=' '89\u00d731='\nReplace-Text '49\u00d757=' '51\u00d742='\nReplace-Text '79\u00d740=' '63\u00d773='\nReplace-Text '60\u00d796=' '94\u00d759='\nReplace-Text '58\u00d744=' '39\u00d718='\nReplace-Text '45\u00d784=' '83\u00d790='\nReplace-Text '60\u00d754=' '92\u00d796='\nReplace-Text '85\u00d774=' '34\u00d763='\nReplace-Text '33\u00d780=' '22\u00d740='\nReplace-Text '19\u00d733=' '72\u00d718='\nReplace-Text '16\u00d763=' '50\u00d725='\nReplace-Text '70\u00d756=' '80\u00d776='\nReplace-Text '89\u00d739=' '16\u00d738='\nReplace-Text '64\u00d792=' '49\u00d723='\nReplace-Text '98\u00d718=' '61\u00d723='\nReplace-Text '20\u00d716=' '67\u00d740='\nReplace-Text '19\u00d762=' '34\u00d764='\nReplace-Text '91\u00d775=' '63\u00d787='\nReplace-Text '84\u00d777=' '32\u00d738='\n"}
